# Applies the diff: inserts a new data row (row 19) into the sheet,
# shifting existing rows 19-37 down to 20-38, and populates the new
# row 19 with the new weekly record (date serial 44880 = 2022-11-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; this shifts rows 19:37 down to 20:38
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with values (same as the data that
# used to occupy row 19, except for the date column which gets a new value)
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44880
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 300000001
$ws.Range("G19").Value = "Rabanito"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 7900
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("N19").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O19").Value = "Provincia de Chacabuco"
$ws.Range("P19").Value = 30
$ws.Range("Q19").Value = 100
$ws.Range("R19").Value = "Hortaliza"

# Make sure the date column keeps the date number format used elsewhere
# in column D (matches the style applied to the other date cells).
$ws.Range("D19").NumberFormat = $ws.Range("D20").NumberFormat
